# Updating the models for NRG, PCSun and Ulmeni
# Shift all timestamps in column A forward by 13 days (new data window),
# and replace the production figures in column B with the refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: shift every timestamp (rows 2..97) forward by 13 days ---
for ($r = 2; $r -le 97; $r++) {
    $oldDate = $ws.Cells.Item($r, 1).Value2()
    $ws.Cells.Item($r, 1).Value = $oldDate + 13
}

# --- Column B: new "Actual Production (MW)" values for rows 2..30 ---
$newB = @(1768,1626,1725,1572,1394,1327,1261,1132,1002,931,865,836,807,752,723,673,604,585,577,549,349,300,292,275,249,221,213,171,140)

$r = 2
foreach ($val in $newB) {
    $ws.Cells.Item($r, 2).Value = $val
    $r++
}

# --- Column B: rows 31..97 are all zero in the refreshed dataset ---
for ($r = 31; $r -le 97; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
}
